# fix : to_excel print 삭제
# Update stack usage counts across the STACK_LIST sheets.

$wb = $excel.ActiveWorkbook

# FRONT_STACK_LIST
$ws = $wb.Worksheets.Item("FRONT_STACK_LIST")
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 3
$ws.Range("B9").Value = 3
$ws.Range("B10").Value = 2

# BACK_STACK_LIST
$ws = $wb.Worksheets.Item("BACK_STACK_LIST")
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 3
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 2

# SECURITY_STACK_LIST
$ws = $wb.Worksheets.Item("SECURITY_STACK_LIST")
$ws.Range("B2").Value = 3

# CLOUD_STACK_LIST
$ws = $wb.Worksheets.Item("CLOUD_STACK_LIST")
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 3
$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 3
$ws.Range("B10").Value = 2
